# Update cryptos price (D) and volume-change (E) columns with latest snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.324.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -7.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.677.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5099"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -13.08%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.02"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06321"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07365"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.681.57"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.10%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5781"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.908.82"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008533"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.65"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -13.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.355.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -7.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.003"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.21%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "186.23"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -9.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.230"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.86%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.85"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.469"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -7.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1170"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.83"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.337"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05810"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.326"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.511"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.507"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.74%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.88%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5939"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.358"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.666"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.094.50"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.887"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8603"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.91"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.833.68"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000114"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.31"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.976"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.80%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05208"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.96%  "
